$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column L ("2021") that duplicates column K ("2020").

# Row 3: empty separator cell with a bottom border (same as K3).
$b3 = $ws.Range("L3").Borders.Item(9)
$b3.Weight = -4138
$b3.Color = 0

# Row 4: year header "2021", bold, right aligned, with bottom border (same as K4).
$ws.Range("L4").Value = 2021
$ws.Range("L4").Font.Bold = $true
$ws.Range("L4").HorizontalAlignment = -4152
$b4 = $ws.Range("L4").Borders.Item(9)
$b4.Weight = -4138
$b4.Color = 0

# Rows 5-10: data values copied from column K.
$ws.Range("L5").Value = 0.86
$ws.Range("L6").Value = 1.07
$ws.Range("L7").Value = 25.27
$ws.Range("L8").Value = 14
$ws.Range("L9").Value = 0.12
$ws.Range("L10").Value = 21.74

# Row 11: total row value, with bottom border (same as K11).
$ws.Range("L11").Value = 9.4600000000000009
$b11 = $ws.Range("L11").Borders.Item(9)
$b11.Weight = -4138
$b11.Color = 0

# Update the active selection.
$ws.Range("N2").Select()
